$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 127
$ws.Range("B3").Value = 127
$ws.Range("B4").Value = 124
$ws.Range("B5").Value = 112
$ws.Range("B6").Value = 107
$ws.Range("B7").Value = 100
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 91
